# g20.3 -> g20.4 : atualização da fonte
# Updates the "Diferença" ranking table to the refreshed source data
# (2010-2000 comparison instead of 2022-2000), re-ordered by value, and
# refreshes the header formatting / page margins to the new template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- refreshed data (rows 2-10), already sorted by the new C-column value ---
$data = @(
    @("Distrito Federal",   "Diferença 2010-2000", -0.08900000000000002, "1º"),
    @("Rio de Janeiro",     "Diferença 2010-2000", -0.08999999999999997, "2º"),
    @("São Paulo",          "Diferença 2010-2000", -0.09100000000000003, "3º"),
    @("Rio Grande do Sul",  "Diferença 2010-2000", -0.093,                "4º"),
    @("Roraima",            "Diferença 2010-2000", -0.09500000000000003, "5º"),
    @("Santa Catarina",     "Diferença 2010-2000", -0.09999999999999998, "6º"),
    @("Sergipe",            "Diferença 2010-2000", -0.138,                "13º"),
    @("Brasil",             "Diferença 2010-2000", -0.12,                 ""),
    @("Nordeste",           "Diferença 2010-2000", -0.11,                 "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

# --- header row formatting refresh: bold font, thin box border, centered + top-aligned ---
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# --- page margins refreshed to Excel defaults (inches -> points: 1pt = 1/72in) ---
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
